# update BSTTCS and TTCS
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text updates (header rows for the two timesheet tables) ---
# First table (row 3 header)
$ws.Range("F3").Value = "Bs_Sapo/ Tỉnh"
$ws.Range("L3").Value = "Bt Duyet Sa/Tối"

# Second table (row 8 header)
$ws.Range("F8").Value = "KTD Sa/Tối"
$ws.Range("H8").Value = "TCT Sa"

# --- Row height adjustments ---
$ws.Rows.Item(3).RowHeight = 39
$ws.Rows.Item(8).RowHeight = 38.25
$ws.Rows.Item(14).RowHeight = 32.1
$ws.Rows.Item(15).RowHeight = 32.1
$ws.Rows.Item(16).RowHeight = 32.1

# Rows 5 and 11 revert to the sheet's default (non-custom) height
$ws.Rows.Item(5).AutoFit() | Out-Null
$ws.Rows.Item(11).AutoFit() | Out-Null

# --- Alignment fix for the "Trừ chỉ tiêu" numeric row (row 10): right -> center ---
$ws.Range("D10:P10").HorizontalAlignment = -4108

# --- Update active selection to match the latest edit location ---
$ws.Range("H8:I8").Select() | Out-Null
